$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P2:P5").NumberFormat = "@"

$ws.Range("P2").Value = "320018616155"
$ws.Range("P3").Value = "320018616166"
$ws.Range("P5").Value = "320018616199"

$ws.Range("P2:P5").Style = "Normal"
